$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.745.76"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.355.30"
$ws.Range("E3").Value = "  +4.72%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.27"
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.658"
$ws.Range("E6").Value = "  +2.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.86"
$ws.Range("E7").Value = "  +14.90%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  +20.11%  "
$ws.Range("E10").Value = "  +2.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "27.35"
$ws.Range("E11").Value = "  +2.48%  "
$ws.Range("E12").Value = "  +2.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.707.49"
$ws.Range("E13").Value = "  +4.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.57"
$ws.Range("E14").Value = "  +11.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.64"
$ws.Range("E15").Value = "  +10.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.878"
$ws.Range("E16").Value = "  +7.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.368.87"
$ws.Range("E17").Value = "  +5.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.623.47"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("E19").Value = "  +4.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.44"
$ws.Range("E20").Value = "  +6.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.60"
$ws.Range("E21").Value = "  +3.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "250.94"
$ws.Range("E22").Value = "  +1.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.83"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.21"
$ws.Range("E26").Value = "  +5.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.24"
$ws.Range("E27").Value = "  -2.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.46"
$ws.Range("E28").Value = "  +4.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.08"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.53"
$ws.Range("E30").Value = "  +6.92%  "
$ws.Range("E31").Value = "  +2.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.131"
$ws.Range("E32").Value = "  +4.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.08"
$ws.Range("E33").Value = "  +3.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0699"
$ws.Range("E34").Value = "  +3.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.10"
$ws.Range("E35").Value = "  +3.62%  "
$ws.Range("E36").Value = "  +2.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.61"
$ws.Range("E37").Value = "  +4.05%  "
$ws.Range("E38").Value = "  +7.39%  "
$ws.Range("E39").Value = "  +5.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.52"
$ws.Range("E40").Value = "  +13.44%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.92"
$ws.Range("E42").Value = "  +2.00%  "
$ws.Range("E43").Value = "  +9.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.25"
$ws.Range("E44").Value = "  +2.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0966"
$ws.Range("E45").Value = "  +3.38%  "
$ws.Range("E46").Value = "  +3.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.45"
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.180"
$ws.Range("E48").Value = "  +12.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.438.53"
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.585.63"
$ws.Range("E50").Value = "  +4.60%  "
$ws.Range("B51").Value = "TerraClassic"
$ws.Range("C51").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000204"
$ws.Range("E51").Value = "  -1.14%  "
